# Add a "Save" column (column H) to the s_vals sheet, matching the diff:
#   - H1 header "Save" styled like the other header cells (bold/border)
#   - H2 = 0, H3 = 1 (plain numeric values, default style)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: set the text first, then clone the header style from G1
# (Copy + PasteSpecial formats only) so H1 ends up with the same bold/
# bordered/centered style used by the rest of row 1.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Data cells: plain numeric values, default (unstyled) formatting.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
